$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows to append (all stored as text, matching existing sheet convention)
$rows = @(
    @("TN330", "Natalie's - Honey Tangerine", "1", "14.57", "14.57"),
    @("AH252", "Natalie's - Orange Juice", "2", "24.50", "49.00"),
    @("TN454", "Natalie's - Orange Mango", "1", "13.38", "13.38"),
    @("TN362", "Natalie's - Orange Pineapple", "1", "13.38", "13.38")
)

$startRow = 17
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}
